$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-18 Monday", "2025-08-19 Tuesday"),
    @("52×94=", "79×82="),
    @("76×58=", "98×63="),
    @("62×82=", "85×58="),
    @("25×42=", "19×27="),
    @("97×47=", "35×61="),
    @("49×57=", "15×59="),
    @("34×36=", "47×97="),
    @("20×30=", "51×49="),
    @("29×46=", "27×30="),
    @("38×14=", "76×84="),
    @("88×61=", "60×84="),
    @("65×89=", "11×82="),
    @("74×68=", "34×86="),
    @("81×98=", "17×51="),
    @("64×96=", "50×94="),
    @("80×71=", "32×85="),
    @("45×79=", "63×62="),
    @("23×93=", "55×96="),
    @("37×92=", "20×93="),
    @("21×33=", "31×94="),
    @("41×97=", "21×47="),
    @("56×79=", "22×63="),
    @("77×65=", "60×54="),
    @("30×38=", "72×82="),
    @("43×88=", "20×44=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
